# De title in StartScene toegevoegd
# Log a new entry on the "week 50" sheet (row 9) for the activity of
# adding the title to StartScene, mirroring the rows above it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("week 50")

# Begintijd / Eindtijd for entry #3 (10:30 - 10:45)
$ws.Range("C9").Value = 0.4375
$ws.Range("D9").Value = 0.44791666666666669

# Activiteiten description (new shared string)
$ws.Range("F9").Value = "De Title in StartScene gezet"

# Leave the selection on the cell that was just edited
$null = $ws.Range("F9").Select()
